$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6852633357048035
$ws.Range("B1").Value = 0.7423092126846313
$ws.Range("C1").Value = 0.8660566806793213
$ws.Range("D1").Value = 1.447826266288757
$ws.Range("E1").Value = 5.014033317565918
